$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203705787658691
$ws.Range("B1").Value = 2.615059852600098
$ws.Range("D1").Value = 2.173031091690063
$ws.Range("E1").Value = 1.168587923049927
